$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 989.89655
$ws.Range("J17").Value = 1068.28
$ws.Range("L17").Value = 3204.84
$ws.Range("N17").Value = -3540.84

$ws.Range("H28").Value = 782.0769
$ws.Range("I28").Value = 594.375
$ws.Range("J28").Value = 1082.4
$ws.Range("K28").Value = 594.375
$ws.Range("L28").Value = 1082.4
$ws.Range("M28").Value = -109.375
$ws.Range("N28").Value = -2052.4

$ws.Range("H113").Value = 2959.7778
$ws.Range("I113").Value = 2823.3333
$ws.Range("J113").Value = 3232.6667
$ws.Range("K113").Value = 2823.3333
$ws.Range("L113").Value = 3232.6667
$ws.Range("M113").Value = 430.6667000000002
$ws.Range("N113").Value = -9740.6667

$ws.Range("H135").Value = 1433.45
$ws.Range("I135").Value = 1440.9286
$ws.Range("J135").Value = 1416
$ws.Range("K135").Value = 12968.3574
$ws.Range("L135").Value = 12744
$ws.Range("M135").Value = -10433.3574
$ws.Range("N135").Value = -17814

$ws.Range("H138").Value = 1962.6289
$ws.Range("I138").Value = 1124.7358
$ws.Range("J138").Value = 2971.9092
$ws.Range("K138").Value = 3374.2074
$ws.Range("L138").Value = 8915.7276
$ws.Range("M138").Value = 1765.7926
$ws.Range("N138").Value = -19195.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2010.3334
$ws.Range("I61").Value = 2224.6
$ws.Range("J61").Value = 1903.2
$ws.Range("K61").Value = 2224.6
$ws.Range("L61").Value = 1903.2
$ws.Range("M61").Value = -2012.6
$ws.Range("N61").Value = -2327.2

$ws.Range("H132").Value = 1765514.9
$ws.Range("I132").Value = 5561.875
$ws.Range("J132").Value = 3209579
$ws.Range("K132").Value = 16685.625
$ws.Range("L132").Value = 9628737
$ws.Range("M132").Value = -14155.625
$ws.Range("N132").Value = -9633797

$ws.Range("H136").Value = 2010.3334
$ws.Range("I136").Value = 2224.6
$ws.Range("J136").Value = 1903.2
$ws.Range("K136").Value = 6673.799999999999
$ws.Range("L136").Value = 5709.6
$ws.Range("M136").Value = -4123.799999999999
$ws.Range("N136").Value = -10809.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H56").Value = 7888.8887
$ws.Range("J56").Value = 7888.8887
$ws.Range("L56").Value = 7888.8887
$ws.Range("N56").Value = -9366.8887

$ws.Range("H82").Value = 11520.846
$ws.Range("I82").Value = 2977.1
$ws.Range("K82").Value = 2977.1
$ws.Range("M82").Value = -2594.1

$ws.Range("H85").Value = 11520.846
$ws.Range("I85").Value = 2977.1
$ws.Range("K85").Value = 2977.1
$ws.Range("M85").Value = -1651.1

$ws.Range("H94").Value = 1779.9286
$ws.Range("I94").Value = 1719.909
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 1719.909
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -1268.909
$ws.Range("N94").Value = -2902

$ws.Range("H99").Value = 1957.1428
$ws.Range("I99").Value = 1851.4286
$ws.Range("J99").Value = 2062.8572
$ws.Range("K99").Value = 1851.4286
$ws.Range("L99").Value = 2062.8572
$ws.Range("M99").Value = -353.4286
$ws.Range("N99").Value = -5058.8572

$ws.Range("H134").Value = 3822.8823
$ws.Range("I134").Value = 3342.1428
$ws.Range("J134").Value = 4159.4
$ws.Range("K134").Value = 10026.4284
$ws.Range("L134").Value = 12478.2
$ws.Range("M134").Value = -7491.428400000001
$ws.Range("N134").Value = -17548.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1763.49
$ws.Range("I31").Value = 1096.1451
$ws.Range("J31").Value = 2852.3157
$ws.Range("K31").Value = 1096.1451
$ws.Range("L31").Value = 2852.3157
$ws.Range("M31").Value = -801.1451
$ws.Range("N31").Value = -3442.3157

$ws.Range("H34").Value = 1763.49
$ws.Range("I34").Value = 1096.1451
$ws.Range("J34").Value = 2852.3157
$ws.Range("K34").Value = 1096.1451
$ws.Range("L34").Value = 2852.3157
$ws.Range("M34").Value = -894.1451
$ws.Range("N34").Value = -3256.3157

$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").ClearContents()
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = 0

$ws.Range("H62").Value = 51600
$ws.Range("I62").Value = 3200
$ws.Range("J62").Value = 100000
$ws.Range("K62").Value = 3200
$ws.Range("L62").Value = 100000
$ws.Range("M62").Value = -2576
$ws.Range("N62").Value = -101248

$ws.Range("H65").Value = 51600
$ws.Range("I65").Value = 3200
$ws.Range("J65").Value = 100000
$ws.Range("K65").Value = 16000
$ws.Range("L65").Value = 500000
$ws.Range("M65").Value = -12880
$ws.Range("N65").Value = -506240

$ws.Range("H132").Value = 2964.7856
$ws.Range("I132").Value = 2224
$ws.Range("J132").Value = 3606.8
$ws.Range("K132").Value = 6672
$ws.Range("L132").Value = 10820.4
$ws.Range("M132").Value = -4142
$ws.Range("N132").Value = -15880.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 16966.666
$ws.Range("I4").Value = 25200
$ws.Range("K4").Value = 75600
$ws.Range("M4").Value = -75488

$ws.Range("H113").Value = 740.58826
$ws.Range("I113").Value = 665.5714
$ws.Range("J113").Value = 793.1
$ws.Range("K113").Value = 1996.7142
$ws.Range("L113").Value = 2379.3
$ws.Range("M113").Value = 173.2857999999999
$ws.Range("N113").Value = -6719.3

$ws.Range("H122").Value = 1712.9697
$ws.Range("J122").Value = 2031.7885
$ws.Range("L122").Value = 18286.0965
$ws.Range("N122").Value = -23186.0965

$ws.Range("H132").Value = 1188.0605
$ws.Range("I132").Value = 1049.7142
$ws.Range("J132").Value = 1290
$ws.Range("K132").Value = 9447.427799999999
$ws.Range("L132").Value = 11610
$ws.Range("M132").Value = -6917.427799999999
$ws.Range("N132").Value = -16670

$ws.Range("H134").Value = 2415.889
$ws.Range("I134").Value = 921.381
$ws.Range("J134").Value = 3723.5833
$ws.Range("K134").Value = 2764.143
$ws.Range("L134").Value = 11170.7499
$ws.Range("M134").Value = 2305.857
$ws.Range("N134").Value = -21310.7499

$ws.Range("H139").Value = 2605.1714
$ws.Range("I139").Value = 1732
$ws.Range("J139").Value = 2750.7
$ws.Range("K139").Value = 5196
$ws.Range("L139").Value = 8252.099999999999
$ws.Range("M139").Value = -56
$ws.Range("N139").Value = -18532.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1191.7273
$ws.Range("I97").Value = 774.875
$ws.Range("J97").Value = 2303.3333
$ws.Range("K97").Value = 774.875
$ws.Range("L97").Value = 2303.3333
$ws.Range("M97").Value = -278.875
$ws.Range("N97").Value = -3295.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 25643810
$ws.Range("I40").Value = 43480624
$ws.Range("J40").Value = 3388.3125
$ws.Range("K40").Value = 43480624
$ws.Range("L40").Value = 3388.3125
$ws.Range("M40").Value = -43480488
$ws.Range("N40").Value = -3660.3125

$ws.Range("H93").Value = 3002.3333
$ws.Range("I93").Value = 3001.5
$ws.Range("K93").Value = 3001.5
$ws.Range("M93").Value = -1753.5

$ws.Range("H100").Value = 2550.4644
$ws.Range("I100").Value = 1804
$ws.Range("K100").Value = 1804
$ws.Range("M100").Value = -1263

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1813.3846
$ws.Range("I132").Value = 1502.5428
$ws.Range("J132").Value = 2453.353
$ws.Range("K132").Value = 4507.6284
$ws.Range("L132").Value = 7360.059
$ws.Range("M132").Value = -1977.6284
$ws.Range("N132").Value = -12420.059
